# Clients.xlsx - "Manque juste les commentaires"
#
# Two placeholder rows (Kosuke Yokono / 364468 and Hana Murata / 584106) that had
# slipped into the bottom of the client table are removed, the remaining twenty
# client rows are put back in their final order, and the last card number
# (123123, Charles Bouvier Dondo) is stored as a plain number instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Drop the two erroneous trailing rows (22: Kosuke Yokono, 23: Hana Murata).
#    Delete from the bottom up so row numbers of the remaining rows do not shift
#    under us while we are still deleting.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

# 2) Stage a copy of the twenty remaining client rows (2-21) well below the table
#    so they can be written back in their new order without a source row ever
#    being overwritten before it has been read.
$ws.Range("A2:D2").Copy()
$ws.Range("A100:D100").PasteSpecial()
$ws.Range("A3:D3").Copy()
$ws.Range("A101:D101").PasteSpecial()
$ws.Range("A4:D4").Copy()
$ws.Range("A102:D102").PasteSpecial()
$ws.Range("A5:D5").Copy()
$ws.Range("A103:D103").PasteSpecial()
$ws.Range("A6:D6").Copy()
$ws.Range("A104:D104").PasteSpecial()
$ws.Range("A7:D7").Copy()
$ws.Range("A105:D105").PasteSpecial()
$ws.Range("A8:D8").Copy()
$ws.Range("A106:D106").PasteSpecial()
$ws.Range("A9:D9").Copy()
$ws.Range("A107:D107").PasteSpecial()
$ws.Range("A10:D10").Copy()
$ws.Range("A108:D108").PasteSpecial()
$ws.Range("A11:D11").Copy()
$ws.Range("A109:D109").PasteSpecial()
$ws.Range("A12:D12").Copy()
$ws.Range("A110:D110").PasteSpecial()
$ws.Range("A13:D13").Copy()
$ws.Range("A111:D111").PasteSpecial()
$ws.Range("A14:D14").Copy()
$ws.Range("A112:D112").PasteSpecial()
$ws.Range("A15:D15").Copy()
$ws.Range("A113:D113").PasteSpecial()
$ws.Range("A16:D16").Copy()
$ws.Range("A114:D114").PasteSpecial()
$ws.Range("A17:D17").Copy()
$ws.Range("A115:D115").PasteSpecial()
$ws.Range("A18:D18").Copy()
$ws.Range("A116:D116").PasteSpecial()
$ws.Range("A19:D19").Copy()
$ws.Range("A117:D117").PasteSpecial()
$ws.Range("A20:D20").Copy()
$ws.Range("A118:D118").PasteSpecial()
$ws.Range("A21:D21").Copy()
$ws.Range("A119:D119").PasteSpecial()

# 3) Re-write rows 2-21 from the staged copies, in their final order.
$ws.Range("A100:D100").Copy()
$ws.Range("A2:D2").PasteSpecial()
$ws.Range("A112:D112").Copy()
$ws.Range("A3:D3").PasteSpecial()
$ws.Range("A104:D104").Copy()
$ws.Range("A4:D4").PasteSpecial()
$ws.Range("A106:D106").Copy()
$ws.Range("A5:D5").PasteSpecial()
$ws.Range("A102:D102").Copy()
$ws.Range("A6:D6").PasteSpecial()
$ws.Range("A109:D109").Copy()
$ws.Range("A7:D7").PasteSpecial()
$ws.Range("A114:D114").Copy()
$ws.Range("A8:D8").PasteSpecial()
$ws.Range("A101:D101").Copy()
$ws.Range("A9:D9").PasteSpecial()
$ws.Range("A117:D117").Copy()
$ws.Range("A10:D10").PasteSpecial()
$ws.Range("A115:D115").Copy()
$ws.Range("A11:D11").PasteSpecial()
$ws.Range("A105:D105").Copy()
$ws.Range("A12:D12").PasteSpecial()
$ws.Range("A110:D110").Copy()
$ws.Range("A13:D13").PasteSpecial()
$ws.Range("A108:D108").Copy()
$ws.Range("A14:D14").PasteSpecial()
$ws.Range("A107:D107").Copy()
$ws.Range("A15:D15").PasteSpecial()
$ws.Range("A113:D113").Copy()
$ws.Range("A16:D16").PasteSpecial()
$ws.Range("A116:D116").Copy()
$ws.Range("A17:D17").PasteSpecial()
$ws.Range("A118:D118").Copy()
$ws.Range("A18:D18").PasteSpecial()
$ws.Range("A103:D103").Copy()
$ws.Range("A19:D19").PasteSpecial()
$ws.Range("A111:D111").Copy()
$ws.Range("A20:D20").PasteSpecial()
$ws.Range("A119:D119").Copy()
$ws.Range("A21:D21").PasteSpecial()

# 4) Clean up the staging area.
$ws.Range("A100:D119").Clear()

# 5) The last card number is now entered as a genuine number, not text.
$ws.Cells.Item(21, 1).Value = 123123

